$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 15:05"
$ws.Range("B4").Value = 1667935
$ws.Range("C4").Value = 1107
$ws.Range("E4").Value = 1122303
$ws.Range("G4").Value = 22
$ws.Range("H4").Value = 98705
$ws.Range("B11").Value = 180020
$ws.Range("C11").Value = 34
$ws.Range("E11").Value = 11351
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 8369
$ws.Range("B14").Value = 132919
$ws.Range("C14").Value = 1496
$ws.Range("D14").Value = 54865
$ws.Range("E14").Value = 74153
$ws.Range("G14").Value = 33
$ws.Range("H14").Value = 3901
$ws.Range("B18").Value = 72560
$ws.Range("C18").Value = 2399
$ws.Range("D18").Value = 43520
$ws.Range("E18").Value = 28650
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 390
$ws.Range("B24").Value = 43714
$ws.Range("C24").Value = 1501
$ws.Range("D24").Value = 9170
$ws.Range("E24").Value = 34521
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 23
$ws.Range("E40").Value = 5486
$ws.Range("G40").Value = 9
$ws.Range("H40").Value = 1185
$ws.Range("B50").Value = 11159
$ws.Range("C50").Value = 67
$ws.Range("D50").Value = 5857
$ws.Range("E50").Value = 5064
$ws.Range("D56").Value = 4320
$ws.Range("E56").Value = 3967
$ws.Range("A91").Value = "El Salvador"
$ws.Range("B91").Value = 1915
$ws.Range("C91").Value = 96
$ws.Range("D91").Value = 594
$ws.Range("E91").Value = 1286
$ws.Range("G91").Value = 2
$ws.Range("H91").Value = 35
$ws.Range("A92").Value = "Estonia"
$ws.Range("B92").Value = 1823
$ws.Range("C92").Value = 2
$ws.Range("D92").Value = 1532
$ws.Range("E92").Value = 227
$ws.Range("G92").Value = 0
$ws.Range("H92").Value = 64
$ws.Range("B103").Value = 1118
$ws.Range("C103").Value = 29
$ws.Range("E103").Value = 435
$ws.Range("B109").Value = 1030
$ws.Range("C109").Value = 15
$ws.Range("D109").Value = 597
$ws.Range("E109").Value = 368
$ws.Range("G109").Value = 2
$ws.Range("H109").Value = 65
$ws.Range("A124").Value = "Sierra Leona"
$ws.Range("B124").Value = 707
$ws.Range("C124").Value = 86
$ws.Range("D124").Value = 241
$ws.Range("E124").Value = 426
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 40
$ws.Range("A125").Value = "Jordania"
$ws.Range("B125").Value = 704
$ws.Range("D125").Value = 470
$ws.Range("E125").Value = 225
$ws.Range("H125").Value = 9
$ws.Range("A126").Value = "San Marino"
$ws.Range("B126").Value = 665
$ws.Range("D126").Value = 266
$ws.Range("E126").Value = 357
$ws.Range("H126").Value = 42
$ws.Range("A127").Value = "Sudan del Sur"
$ws.Range("B127").Value = 655
$ws.Range("D127").Value = 6
$ws.Range("E127").Value = 641
$ws.Range("H127").Value = 8
$ws.Range("A128").Value = "Republica del Chad"
$ws.Range("B128").Value = 648
$ws.Range("D128").Value = 204
$ws.Range("E128").Value = 384
$ws.Range("H128").Value = 60
$ws.Range("D157").Value = 51
$ws.Range("E157").Value = 143
$ws.Range("A173").Value = "Siria"
$ws.Range("B173").Value = 86
$ws.Range("C173").Value = 16
$ws.Range("D173").Value = 41
$ws.Range("E173").Value = 41
$ws.Range("G173").Value = 0
$ws.Range("A174").Value = "Malaui"
$ws.Range("D174").Value = 28
$ws.Range("E174").Value = 50
$ws.Range("G174").Value = 1
$ws.Range("H174").Value = 4
$ws.Range("A175").Value = "Liechtenstein"
$ws.Range("B175").Value = 82
$ws.Range("D175").Value = 55
$ws.Range("E175").Value = 26
$ws.Range("A176").Value = "Comoras"
$ws.Range("B176").Value = 78
$ws.Range("D176").Value = 18
$ws.Range("E176").Value = 59
$ws.Range("H176").Value = 1
$ws.Range("A177").Value = "San Martin (Parte Holandesa)"
$ws.Range("B177").Value = 77
$ws.Range("D177").Value = 59
$ws.Range("E177").Value = 3
$ws.Range("H177").Value = 15
$ws.Range("A178").Value = "Libia"
$ws.Range("B178").Value = 75
$ws.Range("D178").Value = 39
$ws.Range("E178").Value = 33
$ws.Range("H178").Value = 3
$ws.Range("A198").Value = "Nueva Caledonia"
$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0
$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1
$ws.Range("A214").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A215").Value = "San Bartolome"
